$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of an existing category cell so the new row matches
# the rest of the list (font "Times New Roman", size 12, color FF333333 -> style index 1)
$ws.Range("A2").Copy()
$ws.Range("A14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Add the new question category
$ws.Range("A14").Value = "Laboratório"

# Move the active selection to A15, as Excel does after entering data and pressing Enter
$ws.Range("A15").Select()
